$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and re-ordered rows for dogwifhat/Hedera and Stellar/Mantle)

$ws.Range("D2").Value = "68.143.52"
$ws.Range("E2").Value = "  +2.94%  "

$ws.Range("D3").Value = "3.311.45"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.38"
$ws.Range("E5").Value = "  +3.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.44"
$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("D9").Value = "3.306.89"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +0.69%  "

$ws.Range("E11").Value = "  +1.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.12"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +3.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "672.44"
$ws.Range("E14").Value = "  +10.03%  "

$ws.Range("D15").Value = "3.845.91"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "68.244.31"
$ws.Range("E17").Value = "  +3.10%  "

$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").Value = "3.315.29"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.98"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.902"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.75"
$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.12"
$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.91"
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.60"
$ws.Range("E29").Value = "  +5.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.58"
$ws.Range("E30").Value = "  +0.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("E31").Value = "  +2.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "599.69"
$ws.Range("E32").Value = "  +6.71%  "

$ws.Range("D33").Value = "3.975.27"
$ws.Range("E33").Value = "  +3.80%  "

$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("B35").Value = "dogwifhat"
$ws.Range("C35").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.54"
$ws.Range("E35").Value = "  -4.64%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("E36").Value = "  +1.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.99"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.132"
$ws.Range("E39").Value = "  +3.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("E40").Value = "  +5.07%  "

$ws.Range("E41").Value = "  +3.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.71"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("D43").Value = "0.0₃0688"
$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.40"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.337"
$ws.Range("E45").Value = "  +1.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0418"
$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.42"
$ws.Range("E47").Value = "  +12.59%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  +1.63%  "

$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.54"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.77"
$ws.Range("E51").Value = "  +2.18%  "
